$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the email address text (shared string used by A2), replacing the
# old hyperlink's display text with the new one.
$ws.Range("A2").Value = "joshua.lee.hunter@hpe.com"

# Remove the mailto hyperlink that used to be attached to A2.
$ws.Hyperlinks.Delete()

# Move the active selection to A6 (no hyperlink there).
$ws.Range("A6").Select()
